$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Only the rows whose Price (column D) changes are touched here; D45 is left
# completely alone since its price did not change in this update.
# NumberFormat is pinned to Text ("@") before each Price write so Excel's
# smart cell-entry parsing does not silently convert numeric-looking price
# strings (e.g. "0.9978") into actual numbers -- the source data is text.
$dPriceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.527.51"
$ws.Range("E2").Value = "  +2.04%  "

$ws.Range("D3").Value = "1.791.84"
$ws.Range("E3").Value = "  +4.45%  "

$ws.Range("D4").Value = "0.9978"
$ws.Range("E4").Value = "  -0.34%  "

$ws.Range("D5").Value = "314.24"
$ws.Range("E5").Value = "  +1.96%  "

$ws.Range("D6").Value = "0.9984"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").Value = "0.5369"
$ws.Range("E7").Value = "  +12.45%  "

$ws.Range("D8").Value = "0.3787"
$ws.Range("E8").Value = "  +8.88%  "

$ws.Range("D9").Value = "43.05"
$ws.Range("E9").Value = "  +2.66%  "

$ws.Range("D10").Value = "0.07470"
$ws.Range("E10").Value = "  +3.35%  "

$ws.Range("D11").Value = "1.107"
$ws.Range("E11").Value = "  +5.95%  "

$ws.Range("D12").Value = "0.9972"
$ws.Range("E12").Value = "  -0.36%  "

$ws.Range("D13").Value = "20.87"
$ws.Range("E13").Value = "  +5.13%  "

$ws.Range("D14").Value = "6.151"
$ws.Range("E14").Value = "  +5.36%  "

$ws.Range("D15").Value = "1.783.24"
$ws.Range("E15").Value = "  +3.94%  "

$ws.Range("D16").Value = "7.037"
$ws.Range("E16").Value = "  +2.91%  "

$ws.Range("D17").Value = "90.24"
$ws.Range("E17").Value = "  +4.32%  "

$ws.Range("D18").Value = "0.00001062"
$ws.Range("E18").Value = "  +2.38%  "

$ws.Range("D19").Value = "0.06454"
$ws.Range("E19").Value = "  +1.18%  "

$ws.Range("D20").Value = "0.9982"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").Value = "16.92"
$ws.Range("E21").Value = "  +2.55%  "

$ws.Range("D22").Value = "5.938"
$ws.Range("E22").Value = "  +5.77%  "

$ws.Range("D23").Value = "27.521.08"
$ws.Range("E23").Value = "  +1.81%  "

$ws.Range("D24").Value = "11.23"
$ws.Range("E24").Value = "  +4.51%  "

$ws.Range("D25").Value = "2.094"
$ws.Range("E25").Value = "  +0.28%  "

$ws.Range("D26").Value = "156.36"
$ws.Range("E26").Value = "  +3.25%  "

$ws.Range("D27").Value = "20.50"
$ws.Range("E27").Value = "  +2.75%  "

$ws.Range("D28").Value = "2.405"
$ws.Range("E28").Value = "  +15.47%  "

$ws.Range("D29").Value = "1.985.88"
$ws.Range("E29").Value = "  +4.45%  "

$ws.Range("D30").Value = "121.79"
$ws.Range("E30").Value = "  +0.75%  "

$ws.Range("D31").Value = "1.099"
$ws.Range("E31").Value = "  +6.81%  "

$ws.Range("D32").Value = "0.1025"
$ws.Range("E32").Value = "  +12.19%  "

$ws.Range("D33").Value = "5.636"
$ws.Range("E33").Value = "  +5.69%  "

$ws.Range("D34").Value = "3.622"
$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("D35").Value = "0.02273"
$ws.Range("E35").Value = "  +4.37%  "

$ws.Range("D36").Value = "0.06013"
$ws.Range("E36").Value = "  +2.44%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "4.963"
$ws.Range("E37").Value = "  +4.93%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").Value = "11.37"
$ws.Range("E38").Value = "  +3.47%  "

$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2073"
$ws.Range("E39").Value = "  +3.59%  "

$ws.Range("D40").Value = "8.358"
$ws.Range("E40").Value = "  +12.20%  "

$ws.Range("D41").Value = "0.6176"
$ws.Range("E41").Value = "  +2.00%  "

$ws.Range("D42").Value = "1.419"
$ws.Range("E42").Value = "  -3.21%  "

$ws.Range("D43").Value = "0.9974"
$ws.Range("E43").Value = "  -0.27%  "

$ws.Range("D44").Value = "1.145"
$ws.Range("E44").Value = "  +5.36%  "

$ws.Range("E45").Value = "  +3.82%  "

$ws.Range("D46").Value = "0.5822"
$ws.Range("E46").Value = "  +3.45%  "

$ws.Range("D47").Value = "3.634"
$ws.Range("E47").Value = "  +1.86%  "

$ws.Range("D48").Value = "121.51"
$ws.Range("E48").Value = "  +2.07%  "

$ws.Range("D49").Value = "1.909"
$ws.Range("E49").Value = "  +4.22%  "

$ws.Range("D50").Value = "1.132"
$ws.Range("E50").Value = "  +1.92%  "

$ws.Range("D51").Value = "0.06745"
$ws.Range("E51").Value = "  +1.40%  "
